# --------------------------------------------------------------------
# Move the paragraph about Schelling's "unnecessary condition" further
# down the section, and put new text about the tolerance-range finding
# in its former place (keeping the _GoBack bookmark with the new text).
# --------------------------------------------------------------------

$d = $word.ActiveDocument

$oldText = "An interesting finding of our work, showing the value of formalizing verbal reasoning, is the Schelling added an unnecessary condition: that the new neighborhood to which an agent moves most be acceptable, since merely random moves suffice to produce the phenomenon he describes. Furthermore, that unnecessary condition can cause our model to run forever: it is quite possible that there is no acceptable neighborhood for some combination of parameters, so that an attempt to randomly move to one, with no check on the number of attempts, will never terminate."
$newText = "By making the tolerance level a range, rather than a single scalar, we see neighborhoods that are “ragged at the edges”: instead of clean divides when we reach equilibrium, like we got with a single tolerance number, there are a scattering of highly tolerant agents “hanging around” the edges of neighborhoods where they are a distinct minority."

# Locate the paragraph that currently holds the old text by searching for a
# distinctive phrase inside it (robust to the paragraph's numeric index).
$search = $d.Content
$ok = $search.Find.Execute("unnecessary condition can cause our model", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $search.Paragraphs(1)
$targetIndex = $targetPara.Range.Start

# Re-acquire the paragraph via the collection so later references stay
# anchored to document position (not to a "next paragraph" chain, which
# re-resolves live and would otherwise end up pointing at the wrong text
# once earlier text is edited).
$target = $d.Paragraphs(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Start -eq $targetIndex) {
        $target = $cand
    }
}

# Insert a blank spacer paragraph immediately after the target paragraph,
# then a second paragraph (after the spacer) that will hold the relocated
# text. Both inserts happen right after $target, so the bookmark stays
# attached to $target itself.
$void = $target.Range.InsertParagraphAfter()
$spacerIndex = $target.Index + 1
$spacer = $d.Paragraphs($spacerIndex)
$void = $spacer.Range.InsertParagraphAfter()

$relocatedIndex = $spacerIndex + 1
$relocated = $d.Paragraphs($relocatedIndex)
$relocatedRange = $relocated.Range
$void = $relocatedRange.MoveEnd(1, -1)
$relocatedRange.Text = $oldText

# Swap the target paragraph's own text for the new text. Re-fetch the
# paragraph by index first, since the document shifted when we inserted
# the two new paragraphs above.
$target = $d.Paragraphs($target.Index)
$find = $target.Range.Find
$void = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2)
